# Add "Fragment / CDS" as a new available "Type" choice.
#
# This touches:
#  - the "predefined" sheet: a new type-choice value ("Fragment / CDS") is
#    inserted into the C5:C16 list at row 6, pushing the existing entries
#    (old C6:C16) down by one row to C7:C17. The typeChoices named range
#    therefore grows from C5:C16 to C5:C17.
#  - the "Sheet2" (data entry) sheet: one of the example rows (D11) is
#    updated to use the new "Fragment / CDS" type, as a usage example.
#  - the selections (active cell) left behind on both sheets.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Sheet2")
$wsPredefined = $wb.Worksheets.Item("predefined")

# --- "predefined" sheet: insert new type choice at C6, shifting the rest ---
# Capture the existing C6:C16 values (top to bottom) before overwriting,
# then write them back starting at C7 (i.e. shifted down by one row) and
# put the new choice at C6.
$oldTypeValues = @()
for ($r = 6; $r -le 16; $r++) {
    $oldTypeValues += $wsPredefined.Cells.Item($r, 3).Text
}

$wsPredefined.Cells.Item(6, 3).Value = "Fragment / CDS"

for ($i = 0; $i -lt $oldTypeValues.Length; $i++) {
    $wsPredefined.Cells.Item(7 + $i, 3).Value = $oldTypeValues[$i]
}

# Column L (just a blank-space spacer column used by the dropdown helper
# text) already spans rows 5:18 on both sides of the edit, so it is left
# untouched -- only column C grows by the new row 17 entry.

# Update the named range so the dropdown list picks up the new row.
$wb.Names.Item("typeChoices").RefersTo = "=predefined!`$C`$5:`$C`$17"

# --- "Sheet2" sheet: use the new type on the example row (D11) ---
$wsData.Range("D11").Value = "Fragment / CDS"

# --- restore the active cell / selection on each sheet ---
$wsPredefined.Activate() | Out-Null
$wsPredefined.Range("C7").Select() | Out-Null

$wsData.Activate() | Out-Null
$wsData.Range("D16").Select() | Out-Null
